$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 - PANADOL ADVANCE 500 MG 48 TABLETS
$ws.Range("H15").Value = "1:1"

# P15's cell format is numeric (0.00), but the report stores this value as
# literal text ("69.0000") rather than a number - temporarily switch the
# format to Text so the string is preserved verbatim, then restore it.
$ws.Range("P15").NumberFormat = "@"
$ws.Range("P15").Value = "69.0000"
$ws.Range("P15").NumberFormat = "0.00"

$ws.Range("Q15").Value = "0:3"

# Total row
$ws.Range("P20").Value = 404.19999999999999

# Timestamp footer
$ws.Range("A21").Value = "Thursday, 7 August, 2025 10:48 AM"
